$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values to reflect corrected counts
$ws.Range("B7").Value  = "Hearing Aid (53) / Cochlear Implant (17) / None (28)"
$ws.Range("B8").Value  = "Spoken (79) / Total Communication (18) / Cued Speech (1)"
$ws.Range("B9").Value  = "18-100 dB HL (mean (SD): 64 (23))"
$ws.Range("B10").Value = "Yes (16) / No (82)"
$ws.Range("B11").Value = "Female (43) / Male (57)"
$ws.Range("B12").Value = "Yes (36) / No (62)"
$ws.Range("B14").Value = "Unilateral (26) / Bilateral (72)"
$ws.Range("B15").Value = "Yes (34) / No (61)"
$ws.Range("B16").Value = "Full-term (16) / Premature (82)"
$ws.Range("B17").Value = "0-43 services per month (mean (SD): 5 (6))"
$ws.Range("B18").Value = "Sensorineural (62) / Conductive (19) / Mixed (8)"

# Row 19 previously described combined CDI scale; now it's specifically
# "Words and Gestures" and numbers changed
$ws.Range("A19").Value = "Words and Gestures CDI - Words Produced"
$ws.Range("B19").Value = "0-259 words (mean (SD): 33 (53))"

# New row 20 adds the "Words and Sentences" CDI variable
$ws.Range("A20").Value = "Words and Sentences CDI - Words Produced"
$ws.Range("B20").Value = "7-635 words (mean (SD): 148 (184))"
